$p = $ppt.ActivePresentation

# --- Slide 3 ("Solution Architecture Overview"): add a caption textbox below the
#     architecture diagram picture linking to the project's github repo ---
$s3 = $p.Slides.Item(3)

$tb = $s3.Shapes.AddTextbox(1, 48.24, 453.61299212598, 397.4588976378, 21.810944881890)
$tb.Name = "TextBox 2"

$fullText = "github repo: https://github.com/BintangPradana/DetectAnomalyEquipmentAI"
$tb.TextFrame.TextRange.Text = $fullText

$tb.TextFrame.WordWrap = $false
$tb.TextFrame.AutoSize = 1
$tb.Fill.Visible = $false

# Common formatting applied to the whole run first (so it "sticks" across every
# sub-run once the paragraph gets split by the per-run edits below).
$tr = $tb.TextFrame.TextRange
$tr.LanguageID = "en-SG"
$tr.Font.Name = "Calibri"
$tr.Font.NameFarEast = "Calibri"
$tr.Font.NameComplexScript = "Calibri"
$tr.Font.Size = 12

# "github" -> bold
$rGithub = $tb.TextFrame.TextRange.Characters(1, 6)
$rGithub.Font.Bold = $true

# " repo: " -> bold
$rRepo = $tb.TextFrame.TextRange.Characters(7, 7)
$rRepo.Font.Bold = $true

# the URL itself -> hyperlink (not bold)
$rLink = $tb.TextFrame.TextRange.Characters(14, 58)
$rLink.ActionSettings.Item(1).Hyperlink.Address = "https://github.com/BintangPradana/DetectAnomalyEquipmentAI"
